$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H20").Value = 0
$ws.Range("I20").Value = 0
$ws.Range("J20").Value = 0
$ws.Range("K20").Value = 0
$ws.Range("L20").Value = 0
$ws.Range("N20").ClearContents()
$ws.Range("H29").Value = 2852.0667
$ws.Range("I29").Value = 2201.2
$ws.Range("J29").Value = 3177.5
$ws.Range("K29").Value = 6603.599999999999
$ws.Range("L29").Value = 9532.5
$ws.Range("M29").Value = -6322.599999999999
$ws.Range("N29").Value = -10094.5
$ws.Range("H35").Value = 0
$ws.Range("I35").Value = 0
$ws.Range("J35").Value = 0
$ws.Range("K35").Value = 0
$ws.Range("L35").Value = 0
$ws.Range("N35").ClearContents()
$ws.Range("H40").Value = 4682.769
$ws.Range("I40").Value = 10941.2
$ws.Range("J40").Value = 771.25
$ws.Range("K40").Value = 10941.2
$ws.Range("L40").Value = 771.25
$ws.Range("M40").Value = -10766.2
$ws.Range("N40").Value = -1121.25
$ws.Range("H52").Value = 1190
$ws.Range("I52").Value = 950
$ws.Range("J52").Value = 1250
$ws.Range("K52").Value = 2850
$ws.Range("L52").Value = 3750
$ws.Range("M52").Value = -2690
$ws.Range("N52").Value = -4070
$ws.Range("H58").Value = 4619.409
$ws.Range("I58").Value = 169
$ws.Range("J58").Value = 9069.817999999999
$ws.Range("K58").Value = 507
$ws.Range("L58").Value = 27209.454
$ws.Range("M58").Value = -357
$ws.Range("N58").Value = -27509.454
$ws.Range("H64").Value = 3572.7334
$ws.Range("I64").Value = 3399.4285
$ws.Range("J64").Value = 5999
$ws.Range("K64").Value = 3399.4285
$ws.Range("L64").Value = 5999
$ws.Range("M64").Value = -3151.4285
$ws.Range("N64").Value = -6495
$ws.Range("H67").Value = 3572.7334
$ws.Range("I67").Value = 3399.4285
$ws.Range("J67").Value = 5999
$ws.Range("K67").Value = 3399.4285
$ws.Range("L67").Value = 5999
$ws.Range("M67").Value = -2541.4285
$ws.Range("N67").Value = -7715
$ws.Range("H112").Value = 1071.3939
$ws.Range("I112").Value = 799.7143
$ws.Range("J112").Value = 1092.0652
$ws.Range("K112").Value = 2399.1429
$ws.Range("L112").Value = 3276.1956
$ws.Range("M112").Value = -1291.1429
$ws.Range("N112").Value = -5492.1956
$ws.Range("H138").Value = 2200.78
$ws.Range("I138").Value = 1099.1875
$ws.Range("J138").Value = 2410.6072
$ws.Range("K138").Value = 3297.5625
$ws.Range("L138").Value = 7231.821599999999
$ws.Range("M138").Value = 1842.4375
$ws.Range("N138").Value = -17511.8216

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H25").Value = 33692
$ws.Range("I25").Value = 558
$ws.Range("J25").Value = 50259
$ws.Range("K25").Value = 558
$ws.Range("L25").Value = 50259
$ws.Range("N25").Value = -51063
$ws.Range("M25").Value = -156
$ws.Range("H32").Value = 1306791
$ws.Range("I32").Value = 1463398.5
$ws.Range("J32").Value = 116574.2
$ws.Range("K32").Value = 1463398.5
$ws.Range("L32").Value = 116574.2
$ws.Range("M32").Value = -1463111.5
$ws.Range("N32").Value = -117148.2

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H22").Value = 339.2
$ws.Range("I22").Value = 263.14285
$ws.Range("J22").Value = 423.26315
$ws.Range("K22").Value = 263.14285
$ws.Range("L22").Value = 423.26315
$ws.Range("M22").Value = -90.14285000000001
$ws.Range("N22").Value = -769.26315
$ws.Range("H86").Value = 227285.84
$ws.Range("I86").Value = 1407.909
$ws.Range("J86").Value = 779431.9
$ws.Range("K86").Value = 1407.909
$ws.Range("L86").Value = 779431.9
$ws.Range("M86").Value = -284.9090000000001
$ws.Range("N86").Value = -781677.9
$ws.Range("H89").Value = 227285.84
$ws.Range("I89").Value = 1407.909
$ws.Range("J89").Value = 779431.9
$ws.Range("K89").Value = 7039.545
$ws.Range("L89").Value = 3897159.5
$ws.Range("M89").Value = -1423.545
$ws.Range("N89").Value = -3908391.5
$ws.Range("H99").Value = 2544.4375
$ws.Range("I99").Value = 2733.3333
$ws.Range("J99").Value = 1977.75
$ws.Range("K99").Value = 2733.3333
$ws.Range("L99").Value = 1977.75
$ws.Range("M99").Value = -1235.3333
$ws.Range("N99").Value = -4973.75

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value = 4252.4443
$ws.Range("I58").Value = 1528.8572
$ws.Range("J58").Value = 5985.636
$ws.Range("K58").Value = 1528.8572
$ws.Range("L58").Value = 5985.636
$ws.Range("M58").Value = -1325.8572
$ws.Range("N58").Value = -6391.636
$ws.Range("H60").Value = 6531
$ws.Range("I60").Value = 4237.2
$ws.Range("J60").Value = 18000
$ws.Range("K60").Value = 4237.2
$ws.Range("L60").Value = 18000
$ws.Range("M60").Value = -3726.2
$ws.Range("N60").Value = -19022
$ws.Range("H132").Value = 1642.1724
$ws.Range("I132").Value = 1337.5416
$ws.Range("J132").Value = 3104.4
$ws.Range("K132").Value = 4012.6248
$ws.Range("L132").Value = 9313.200000000001
$ws.Range("M132").Value = -1482.6248
$ws.Range("N132").Value = -14373.2
$ws.Range("H136").Value = 4252.4443
$ws.Range("I136").Value = 1528.8572
$ws.Range("J136").Value = 5985.636
$ws.Range("K136").Value = 4586.571599999999
$ws.Range("L136").Value = 17956.908
$ws.Range("M136").Value = -2036.571599999999
$ws.Range("N136").Value = -23056.908

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H40").Value = 301.55
$ws.Range("I40").Value = 222.63637
$ws.Range("J40").Value = 398
$ws.Range("K40").Value = 890.54548
$ws.Range("L40").Value = 1592
$ws.Range("M40").Value = -821.54548
$ws.Range("N40").Value = -1730
$ws.Range("H68").Value = 2180.3257
$ws.Range("I68").Value = 4549.64
$ws.Range("J68").Value = 1209.295
$ws.Range("K68").Value = 13648.92
$ws.Range("L68").Value = 3627.885
$ws.Range("M68").Value = -12837.92
$ws.Range("N68").Value = -5249.885
$ws.Range("H71").Value = 2180.3257
$ws.Range("I71").Value = 4549.64
$ws.Range("J71").Value = 1209.295
$ws.Range("K71").Value = 40946.76
$ws.Range("L71").Value = 10883.655
$ws.Range("M71").Value = -36890.76
$ws.Range("N71").Value = -18995.655
$ws.Range("H113").Value = 1711.5
$ws.Range("I113").Value = 399
$ws.Range("J113").Value = 2149
$ws.Range("K113").Value = 1197
$ws.Range("L113").Value = 6447
$ws.Range("M113").Value = 973
$ws.Range("N113").Value = -10787
$ws.Range("H131").Value = 1010.66156
$ws.Range("I131").Value = 671.8
$ws.Range("J131").Value = 1038.9
$ws.Range("K131").Value = 2015.4
$ws.Range("L131").Value = 3116.7
$ws.Range("M131").Value = 3024.6
$ws.Range("N131").Value = -13196.7
$ws.Range("H132").Value = 1207.3478
$ws.Range("I132").Value = 985.6
$ws.Range("J132").Value = 1471.3334
$ws.Range("K132").Value = 8870.4
$ws.Range("L132").Value = 13242.0006
$ws.Range("M132").Value = -6340.4
$ws.Range("N132").Value = -18302.0006

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H81").Value = 33181
$ws.Range("I81").Value = 0
$ws.Range("J81").Value = 33181
$ws.Range("K81").Value = 0
$ws.Range("L81").Value = 33181
$ws.Range("N81").Value = -35177
$ws.Range("H84").Value = 33181
$ws.Range("I84").Value = 0
$ws.Range("J84").Value = 33181
$ws.Range("K84").Value = 0
$ws.Range("L84").Value = 99543
$ws.Range("N84").Value = -109527

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H34").Value = 18333.666
$ws.Range("I34").Value = 18333.666
$ws.Range("J34").Value = 0
$ws.Range("K34").Value = 18333.666
$ws.Range("L34").Value = 0
$ws.Range("M34").Value = -18130.666
$ws.Range("H62").Value = 3770.2
$ws.Range("I62").Value = 3457
$ws.Range("J62").Value = 4501
$ws.Range("K62").Value = 3457
$ws.Range("L62").Value = 4501
$ws.Range("M62").Value = -2833
$ws.Range("H65").Value = 3770.2
$ws.Range("I65").Value = 3457
$ws.Range("J65").Value = 4501
$ws.Range("K65").Value = 17285
$ws.Range("L65").Value = 22505
$ws.Range("M65").Value = -14165
$ws.Range("H81").Value = 1417.04
$ws.Range("I81").Value = 1113.1111
$ws.Range("J81").Value = 2198.5715
$ws.Range("K81").Value = 2226.2222
$ws.Range("L81").Value = 4397.143
$ws.Range("M81").Value = -1165.2222
$ws.Range("H84").Value = 1417.04
$ws.Range("I84").Value = 1113.1111
$ws.Range("J84").Value = 2198.5715
$ws.Range("K84").Value = 11131.111
$ws.Range("L84").Value = 21985.715
$ws.Range("M84").Value = -5827.111000000001
